$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate Frontend / Frontend Contingency values for rows 5-8
$ws.Range("F5").Value = 5
$ws.Range("G5").Value = 6.5

$ws.Range("F6").Value = 7
$ws.Range("G6").Value = 9.1

$ws.Range("F7").Value = 4000
$ws.Range("G7").Value = 5200

$ws.Range("F8").Value = 5000
$ws.Range("G8").Value = 6500
